$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Files Tab" query text in B4: the two coalesce() defaults that
# used to fall back to the literal 'Not specified in data' now fall back to
# an empty string, matching the merged CDS_Regression query used elsewhere.
$cell = $ws.Range("B4")
$sql = $cell.Value2
$sql = $sql -replace "coalesce\(p\.participant_id, 'Not specified in data'\) as ``Participant ID``,", "coalesce(p.participant_id, '') as ``Participant ID``,"
$sql = $sql -replace "coalesce\(samp\.sample_id, 'Not specified in data'\) as ``Sample ID``,", "coalesce(samp.sample_id, '') as ``Sample ID``,"
$cell.Value = $sql

# Move the active selection from B2 to B3, as recorded in the saved view state.
$ws.Range("B3").Select()
